# Update Metadata sheet: Last Updated timestamp
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:19 PM"

# Update Top Gainers sheet: refreshed ranking data (rows 42-75)
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsGainers.Cells.Item(42, 2).Value = "HITECHGEAR"
$wsGainers.Cells.Item(42, 3).Value = 4.8651
$wsGainers.Cells.Item(42, 4).Value = 2.1287
$wsGainers.Cells.Item(42, 5).Value = 10.9905
$wsGainers.Cells.Item(43, 2).Value = "INDOTHAI"
$wsGainers.Cells.Item(43, 3).Value = 4.8064
$wsGainers.Cells.Item(43, 4).Value = 4.5349
$wsGainers.Cells.Item(43, 5).Value = 43.748
$wsGainers.Cells.Item(44, 2).Value = "SANDUMA"
$wsGainers.Cells.Item(44, 3).Value = 4.593
$wsGainers.Cells.Item(44, 4).Value = 2.1405
$wsGainers.Cells.Item(44, 5).Value = 30.2813
$wsGainers.Cells.Item(45, 2).Value = "LLOYDSENT"
$wsGainers.Cells.Item(45, 3).Value = 4.5646
$wsGainers.Cells.Item(45, 4).Value = 1.8339
$wsGainers.Cells.Item(45, 5).Value = 11.234
$wsGainers.Cells.Item(46, 2).Value = "STAR"
$wsGainers.Cells.Item(46, 3).Value = 4.5025
$wsGainers.Cells.Item(46, 4).Value = 4.4319
$wsGainers.Cells.Item(46, 5).Value = 3.662
$wsGainers.Cells.Item(47, 2).Value = "RECLTD"
$wsGainers.Cells.Item(47, 3).Value = 4.4992
$wsGainers.Cells.Item(47, 4).Value = 3.4756
$wsGainers.Cells.Item(47, 5).Value = 3.4062
$wsGainers.Cells.Item(48, 2).Value = "NBCC"
$wsGainers.Cells.Item(48, 3).Value = 4.4511
$wsGainers.Cells.Item(48, 4).Value = 3.1605
$wsGainers.Cells.Item(48, 5).Value = 7.6018
$wsGainers.Cells.Item(49, 2).Value = "GPPL"
$wsGainers.Cells.Item(49, 3).Value = 4.4154
$wsGainers.Cells.Item(49, 4).Value = 3.4073
$wsGainers.Cells.Item(49, 5).Value = 5.0497
$wsGainers.Cells.Item(50, 2).Value = "BIL"
$wsGainers.Cells.Item(50, 3).Value = 4.3654
$wsGainers.Cells.Item(50, 4).Value = 9.1222
$wsGainers.Cells.Item(50, 5).Value = -0.3203
$wsGainers.Cells.Item(51, 2).Value = "HUDCO"
$wsGainers.Cells.Item(51, 3).Value = 4.3201
$wsGainers.Cells.Item(51, 4).Value = 3.8924
$wsGainers.Cells.Item(51, 5).Value = 5.3884
$wsGainers.Cells.Item(52, 2).Value = "SGMART"
$wsGainers.Cells.Item(52, 3).Value = 4.2736
$wsGainers.Cells.Item(52, 4).Value = 8.2589
$wsGainers.Cells.Item(52, 5).Value = 2.5381
$wsGainers.Cells.Item(53, 2).Value = "MRPL"
$wsGainers.Cells.Item(53, 3).Value = 4.2642
$wsGainers.Cells.Item(53, 4).Value = 9.7103
$wsGainers.Cells.Item(53, 5).Value = 20.0542
$wsGainers.Cells.Item(54, 2).Value = "JKIL"
$wsGainers.Cells.Item(54, 3).Value = 4.1372
$wsGainers.Cells.Item(54, 4).Value = 2.9463
$wsGainers.Cells.Item(54, 5).Value = 1.7584
$wsGainers.Cells.Item(55, 2).Value = "SAMBHV"
$wsGainers.Cells.Item(55, 3).Value = 4.1349
$wsGainers.Cells.Item(55, 4).Value = 2.624
$wsGainers.Cells.Item(55, 5).Value = 5.167
$wsGainers.Cells.Item(56, 2).Value = "SAPPHIRE"
$wsGainers.Cells.Item(56, 3).Value = 4.1265
$wsGainers.Cells.Item(56, 4).Value = 1.7633
$wsGainers.Cells.Item(56, 5).Value = -0.7999
$wsGainers.Cells.Item(57, 2).Value = "PVRINOX"
$wsGainers.Cells.Item(57, 3).Value = 4.1118
$wsGainers.Cells.Item(57, 4).Value = 6.2102
$wsGainers.Cells.Item(57, 5).Value = 14.707
$wsGainers.Cells.Item(58, 2).Value = "KERNEX"
$wsGainers.Cells.Item(58, 3).Value = 4.0782
$wsGainers.Cells.Item(58, 4).Value = 7.542
$wsGainers.Cells.Item(58, 5).Value = 27.2033
$wsGainers.Cells.Item(59, 2).Value = "SUNFLAG"
$wsGainers.Cells.Item(59, 3).Value = 3.997
$wsGainers.Cells.Item(59, 4).Value = 4.333
$wsGainers.Cells.Item(59, 5).Value = 4.6312
$wsGainers.Cells.Item(60, 2).Value = "CMSINFO"
$wsGainers.Cells.Item(60, 3).Value = 3.9096
$wsGainers.Cells.Item(60, 4).Value = 2.6872
$wsGainers.Cells.Item(60, 5).Value = 2.8935
$wsGainers.Cells.Item(61, 2).Value = "GMBREW"
$wsGainers.Cells.Item(61, 3).Value = 3.8999
$wsGainers.Cells.Item(61, 4).Value = -0.53
$wsGainers.Cells.Item(61, 5).Value = 79.029
$wsGainers.Cells.Item(62, 2).Value = "GREENLAM"
$wsGainers.Cells.Item(62, 3).Value = 3.8946
$wsGainers.Cells.Item(62, 4).Value = 3.5858
$wsGainers.Cells.Item(62, 5).Value = 10.721
$wsGainers.Cells.Item(63, 2).Value = "APARINDS"
$wsGainers.Cells.Item(63, 3).Value = 3.8924
$wsGainers.Cells.Item(63, 4).Value = 8.3414
$wsGainers.Cells.Item(63, 5).Value = 15.5876
$wsGainers.Cells.Item(67, 2).Value = "NPST"
$wsGainers.Cells.Item(67, 3).Value = 3.7841
$wsGainers.Cells.Item(67, 4).Value = -2.0689
$wsGainers.Cells.Item(67, 5).Value = -3.5677
$wsGainers.Cells.Item(68, 2).Value = "DCW"
$wsGainers.Cells.Item(68, 3).Value = 3.7544
$wsGainers.Cells.Item(68, 4).Value = 2.3219
$wsGainers.Cells.Item(68, 5).Value = -3.9753
$wsGainers.Cells.Item(69, 2).Value = "RHETAN"
$wsGainers.Cells.Item(69, 3).Value = 3.754
$wsGainers.Cells.Item(69, 4).Value = 4.178
$wsGainers.Cells.Item(69, 5).Value = 6.549
$wsGainers.Cells.Item(70, 2).Value = "HINDPETRO"
$wsGainers.Cells.Item(70, 3).Value = 3.6935
$wsGainers.Cells.Item(70, 4).Value = 6.9335
$wsGainers.Cells.Item(70, 5).Value = 5.7397
$wsGainers.Cells.Item(71, 2).Value = "BHARTIHEXA"
$wsGainers.Cells.Item(71, 3).Value = 3.6718
$wsGainers.Cells.Item(71, 4).Value = 7.0877
$wsGainers.Cells.Item(71, 5).Value = 15.3332
$wsGainers.Cells.Item(72, 2).Value = "HLEGLAS"
$wsGainers.Cells.Item(72, 3).Value = 3.659
$wsGainers.Cells.Item(72, 4).Value = 8.1155
$wsGainers.Cells.Item(72, 5).Value = 27.1239
$wsGainers.Cells.Item(73, 2).Value = "RHIM"
$wsGainers.Cells.Item(73, 3).Value = 3.6544
$wsGainers.Cells.Item(73, 4).Value = 3.2276
$wsGainers.Cells.Item(73, 5).Value = 5.1826
$wsGainers.Cells.Item(74, 2).Value = "SHK"
$wsGainers.Cells.Item(74, 3).Value = 3.6347
$wsGainers.Cells.Item(74, 4).Value = 2.388
$wsGainers.Cells.Item(74, 5).Value = -1.932
$wsGainers.Cells.Item(75, 2).Value = "BCLIND"
$wsGainers.Cells.Item(75, 3).Value = 3.6271
$wsGainers.Cells.Item(75, 4).Value = 2.2945
$wsGainers.Cells.Item(75, 5).Value = 0.1728

# Update 1 Month Performance sheet: refreshed values
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Cells.Item(6, 3).Value = 66.1892
$wsPerf.Cells.Item(36, 2).Value = "MINDTECK"
$wsPerf.Cells.Item(36, 3).Value = 26.9415
$wsPerf.Cells.Item(37, 2).Value = "BHARATWIRE"
$wsPerf.Cells.Item(37, 3).Value = 26.5276
$wsPerf.Cells.Item(38, 2).Value = "HATSUN"
$wsPerf.Cells.Item(38, 3).Value = 26.492
$wsPerf.Cells.Item(39, 2).Value = "INDORAMA"
$wsPerf.Cells.Item(39, 3).Value = 26.4516
$wsPerf.Cells.Item(40, 2).Value = "IFBIND"
$wsPerf.Cells.Item(40, 3).Value = 26.161
